# Weekly update: a new week's record is inserted at the top of the data
# block (row 8, right after the first 6 already-sorted rows), pushing all
# subsequent rows (old 8..41) down by one into rows 9..42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; Excel shifts rows 8:41 down to 9:42.
$ws.Rows.Item(8).Insert()

# Populate the newly-inserted row 8 with this week's record.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44462
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108007
$ws.Range("J8").Value = "Coco"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("Q8").Value = "$/malla 20 unidades"
$ws.Range("R8").Value = "Perú"
$ws.Range("S8").Value = 1200
$ws.Range("T8").Value = 20
